function Set-Row {
    param($Ws, $Rng, $Values)
    $n = $Values.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $Values[$i]
    }
    $Ws.Range($Rng).Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-Row $ws "B2:K2" @(0.9192316907193552, -1.148324593739046, -1.376360942196903, 0.5240065820382563, -0.0027599659219745, 0.09815006083404601, -0.1708892827759409, 0.665940954683469, 0.5666693334229809, 0.2125078032274952)
Set-Row $ws "B3:K3" @(-0.740092140117276, -1.124737898310509, 0.442100528215121, 0.07233208345931114, 0.2040572213477698, -0.1124726366172855, 0.7323796385135606, 0.6413512920338174, 0.2828210951284739, 0.4305191322240596)
Set-Row $ws "B4:K4" @(-1.230237110159324, 0.3159952442497632, 0.2405801520925428, 0.2647213139516521, -0.1013397972426526, 0.7952423898404911, 0.6987028198275607, 0.3273262721328529, 0.4815930032504782, 0.6412615601838532)
Set-Row $ws "B5:K5" @(-0.2884453755979755, 0.1763767747378494, 0.5562686220088688, -0.1390466733285862, 0.78665836976703, 0.7824438974249289, 0.3621534170608979, 0.5081573863027752, 0.686760400930452, 0.1781580448571292)
Set-Row $ws "B6:K6" @(0.08873594589893813, 0.5038494199792003, -0.0804002763435806, 0.7915875151041314, 0.7713506143760337, 0.3749201162859844, 0.5168869453501342, 0.6892623450793038, 0.1844339823288103, 0.4621062954844631)
Set-Row $ws "B7:K7" @(0.526277399612209, -0.07779673678113191, 0.7695198294501161, 0.7653686340716788, 0.3695768661035733, 0.506518411979768, 0.68078502295568, 0.1766616740219995, 0.4534823219514945, 0.2839431369332225)
Set-Row $ws "B8:K8" @(-0.05637216532391182, 0.8896976782493284, 0.6761742896578956, 0.3398498622549955, 0.5151445320096781, 0.6548448495302448, 0.153883110993772, 0.4398642868028766, 0.2654223397480467, 0.570669944985061)
Set-Row $ws "B9:K9" @(0.847377045928939, 0.6469698158021624, 0.3368292624500743, 0.4962832483981977, 0.6345137184650405, 0.1393529555595242, 0.4232232413106087, 0.2478384943192965, 0.5541960614550182, -0.05728328644410208)
Set-Row $ws "B10:K10" @(0.9871675564200725, 0.4122003242340114, 0.3051899620851986, 0.6592076310517737, 0.1487323591158202, 0.3684124426992176, 0.2289444034306267, 0.5393320606399725, -0.0859949970734728, 0.6071339948549791)
Set-Row $ws "B11:K11" @(0.862895196224262, 0.352300664297557, 0.4189244002609654, 0.184677440181683, 0.3663616852596248, 0.1523761639945965, 0.511932666264689, -0.1122840472711982, 0.5637367041416466, 0.2970525035592049)
Set-Row $ws "B12:J12" @(0.6646262512210954, 0.5522131399964898, 0.005599018365491398, 0.394555924030192, 0.18840147518699, 0.4743910960604755, -0.1178263863585594, 0.5679120330803951, 0.2858677898194339)
Set-Row $ws "B13:I13" @(0.787803631104331, 0.09027775923980097, 0.2489555573964748, 0.1983552180462326, 0.4897562657600204, -0.1566747213159825, 0.5520662240532093, 0.2775335613519331)
Set-Row $ws "B14:H14" @(0.4013019457211838, 0.36604433180767, 0.03589107659666579, 0.518715216225222, -0.1209318488610789, 0.5207385776695821, 0.2743085116504074)
Set-Row $ws "B15:G15" @(0.6128695092117844, 0.0535469441345553, 0.4240929771142275, -0.0887144606125988, 0.5331267034972994, 0.2534447081011285)
Set-Row $ws "B16:F16" @(0.2932233035507672, 0.509117220583441, -0.2051990389706129, 0.5482319751491519, 0.2766837437271186)
Set-Row $ws "B17:E17" @(0.6739772976175282, -0.1916617667226967, 0.4852590561591889, 0.2867219094086165)
Set-Row $ws "B18:D18" @(0.06218727514271133, 0.5845771063412253, 0.1751453671933744)
Set-Row $ws "B19:C19" @(0.6286367975806744, 0.1965658720679752)
Set-Row $ws "B20:B20" @(0.4328090033804217)

Write-Output "Applied naive component forecaster fix."
